# Natmi following Dr Hou advice
# Rebuild the LR-pair table so that the "Sending cluster" dimension now
# includes ECs (in addition to the existing FAPs and sCs), and refresh all
# of the computed numeric columns for the full cross-product of
# Sending cluster x Target cluster (ECs, FAPs, sCs), keeping
# Ligand symbol = Il4 and Receptor symbol = Il13ra1 fixed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# then the 16 numeric metric columns (E..T).
$rows = @(
    @("ECs", "Il4", "Il13ra1", "ECs", 2, 0.6666666666666666, 0.6580480000000001, 1.974144, 0.1830904640197835, 0.1830904640197835, 3, 1, 25.41813, 76.25439, 0.6391016201831307, 0.6391016201831305, 16.72634961024, 150.53714649216, 0.1170134121951248, 0.1170134121951248),
    @("ECs", "Il4", "Il13ra1", "FAPs", 2, 0.6666666666666666, 0.6580480000000001, 1.974144, 0.1830904640197835, 0.1830904640197835, 3, 1, 10.37350166666667, 31.120505, 0.2608264936145607, 0.2608264936145606, 6.826262024746669, 61.43635822272001, 0.047754843744543, 0.04775484374454299),
    @("ECs", "Il4", "Il13ra1", "sCs", 2, 0.6666666666666666, 0.6580480000000001, 1.974144, 0.1830904640197835, 0.1830904640197835, 3, 1, 3.980024666666667, 11.940074, 0.1000718862023088, 0.1000718862023088, 2.619047271850667, 23.571425446656, 0.01832220808011568, 0.01832220808011568),
    @("FAPs", "Il4", "Il13ra1", "ECs", 3, 1, 1.115367, 3.346101, 0.3103315587652478, 0.3103315587652479, 3, 1, 25.41813, 76.25439, 0.6391016201831307, 0.6391016201831305, 28.35054340371, 255.15489063339, 0.1983334020008263, 0.1983334020008263),
    @("FAPs", "Il4", "Il13ra1", "FAPs", 3, 1, 1.115367, 3.346101, 0.3103315587652478, 0.3103315587652479, 3, 1, 10.37350166666667, 31.120505, 0.2608264936145607, 0.2608264936145606, 11.570261433445, 104.132352901005, 0.08094269233068058, 0.08094269233068056),
    @("FAPs", "Il4", "Il13ra1", "sCs", 3, 1, 1.115367, 3.346101, 0.3103315587652478, 0.3103315587652479, 3, 1, 3.980024666666667, 11.940074, 0.1000718862023088, 0.1000718862023088, 4.439188172386, 39.952693551474, 0.03105546443374098, 0.03105546443374098),
    @("sCs", "Il4", "Il13ra1", "ECs", 3, 1, 1.820699, 5.462097, 0.5065779772149687, 0.5065779772149687, 3, 1, 25.41813, 76.25439, 0.6391016201831307, 0.6391016201831305, 46.27876387287, 416.50887485583, 0.3237548059871795, 0.3237548059871795),
    @("sCs", "Il4", "Il13ra1", "FAPs", 3, 1, 1.820699, 5.462097, 0.5065779772149687, 0.5065779772149687, 3, 1, 10.37350166666667, 31.120505, 0.2608264936145607, 0.2608264936145606, 18.88702411099834, 169.983216998985, 0.1321289575393371, 0.1321289575393371),
    @("sCs", "Il4", "Il13ra1", "sCs", 3, 1, 1.820699, 5.462097, 0.5065779772149687, 0.5065779772149687, 3, 1, 3.980024666666667, 11.940074, 0.1000718862023088, 0.1000718862023088, 7.246426930575334, 65.21784237517799, 0.05069421368845212, 0.05069421368845211)
)

$rowIndex = 2
foreach ($row in $rows) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex = $colIndex + 1
    }
    $rowIndex = $rowIndex + 1
}
